# Minor wording changes to save space
#
# 1) Shape "TextBox 5"  (// x.f = y)       ->  // x.f := y
# 2) Shape "TextBox 7"  (// y = x.f)       ->  // y := x.f
# 3) Shape "TextBox 9"  (Initialize(...))  -> merge two adjacent runs of spaces
# 4) Shape "TextBox 10" (Eq(...))          -> add trailing comment, widen shape

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) "// x.f = y"  ->  "// x.f := y"
#    (split the trailing " = y" run into " " + ":= " + "y")
# ---------------------------------------------------------------------------
$shp1 = $s.Shapes.Item("TextBox 5")
$tr1 = $shp1.TextFrame.TextRange
$mid1 = $tr1.Characters(8, 2)   # "= "
$mid1.Text = ":= "
# The extra character re-triggers this autofit textbox's height calc
# (re-wraps to an extra line); the real deck keeps the original box
# height, so restore it explicitly.
$shp1.Height = 55.73905563354492

# ---------------------------------------------------------------------------
# 2) "// y = x.f"  ->  "// y := x.f"
#    (split the leading "// y = " run into "// y " + ":= ")
# ---------------------------------------------------------------------------
$shp2 = $s.Shapes.Item("TextBox 7")
$tr2 = $shp2.TextFrame.TextRange
$mid2 = $tr2.Characters(6, 2)   # "= "
$mid2.Text = ":= "

# ---------------------------------------------------------------------------
# 3) Merge the "                 " run and the "linear " run into one run
# ---------------------------------------------------------------------------
$shp3 = $s.Shapes.Item("TextBox 9")
$tr3 = $shp3.TextFrame.TextRange
$combined3 = $tr3.Characters(32, 24)  # "                 linear "
$combined3.Text = "                 linear "

# ---------------------------------------------------------------------------
# 4) "Eq(linear tid:Tid, x:idx, y:idx)" -> "...y:idx) // x == y"
#    and widen/reposition the textbox so the longer line still fits.
# ---------------------------------------------------------------------------
$shp4 = $s.Shapes.Item("TextBox 10")
$tr4 = $shp4.TextFrame.TextRange
$tail4 = $tr4.Characters(32, 1)  # ")"
$tail4.Text = ") // x == y"

$shp4.Left = 166.17489624023438
$shp4.Width = 132.70851135253906
